$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B of the "Rules" sheet holds the rule name for the last
# rule (previously "R40"); the commit changes it to the text "1".
# A leading apostrophe forces Excel to store this as text instead of a
# number, matching the <c t="s"> (shared-string / text) cell in the diff.
$ws.Range("B11").Value = "'1"
